$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 381, shifting existing rows 381:401 down to 384:404
$ws.Range("A381:R383").Insert()

# New data block for date 44585 (rows 381-383)
$ws.Range("A381").Value = 3
$ws.Range("B381").Value = "Femacal de La Calera"
$ws.Range("C381").Value = "Coquimbo"
$ws.Range("D381").Value = 44585
$ws.Range("E381").Value = 5
$ws.Range("F381").Value = 100112028
$ws.Range("G381").Value = "Sandia"
$ws.Range("H381").Value = "Sin especificar"
$ws.Range("I381").Value = "Extra"
$ws.Range("J381").Value = 310
$ws.Range("K381").Value = 2500
$ws.Range("L381").Value = 3000
$ws.Range("M381").Value = 2758
$ws.Range("N381").Value = "`$/unidad"
$ws.Range("O381").Value = "Paine"
$ws.Range("P381").Value = 2758
$ws.Range("Q381").Value = 1
$ws.Range("R381").Value = "Hortaliza"

$ws.Range("A382").Value = 3
$ws.Range("B382").Value = "Femacal de La Calera"
$ws.Range("C382").Value = "Coquimbo"
$ws.Range("D382").Value = 44585
$ws.Range("E382").Value = 5
$ws.Range("F382").Value = 100112028
$ws.Range("G382").Value = "Sandia"
$ws.Range("H382").Value = "Sin especificar"
$ws.Range("I382").Value = "Primera"
$ws.Range("J382").Value = 160
$ws.Range("K382").Value = 2000
$ws.Range("L382").Value = 2000
$ws.Range("M382").Value = 2000
$ws.Range("N382").Value = "`$/unidad"
$ws.Range("O382").Value = "Paine"
$ws.Range("P382").Value = 2000
$ws.Range("Q382").Value = 1
$ws.Range("R382").Value = "Hortaliza"

$ws.Range("A383").Value = 3
$ws.Range("B383").Value = "Femacal de La Calera"
$ws.Range("C383").Value = "Coquimbo"
$ws.Range("D383").Value = 44585
$ws.Range("E383").Value = 5
$ws.Range("F383").Value = 100112028
$ws.Range("G383").Value = "Sandia"
$ws.Range("H383").Value = "Sin especificar"
$ws.Range("I383").Value = "Segunda"
$ws.Range("J383").Value = 160
$ws.Range("K383").Value = 1500
$ws.Range("L383").Value = 1500
$ws.Range("M383").Value = 1500
$ws.Range("N383").Value = "`$/unidad"
$ws.Range("O383").Value = "Paine"
$ws.Range("P383").Value = 1500
$ws.Range("Q383").Value = 1
$ws.Range("R383").Value = "Hortaliza"

# Match the date cell format (style) used by the rest of column D
$ws.Range("D384").Copy()
$ws.Range("D381:D383").PasteSpecial(-4122)
$excel.CutCopyMode = 0
